$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the two mistyped column headers.
$ws.Range("B1").Value = "comments"
$ws.Range("C1").Value = "moderator_rating"

# Remove the extraneous last data row (video 7yzImCTv7A0), shrinking the
# sheet's used range from A1:D10 down to A1:D9.
$ws.Rows.Item(10).Delete()
